$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ClassName (column D) for rows 37-41 to include the new "Daily." segment,
# reflecting that these are now the "Daily" accrual leave-balance test cases.
$ws.Range("D37").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalanceTests.Daily.LeaveBalance_48EmployeeCreation_1_25"
$ws.Range("D38").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalanceTests.Daily.LeaveBalance_48EmployeeCreation_26_41"
$ws.Range("D39").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalanceTests.Daily.LeaveBalance_48EmployeeCreation_42_82"
$ws.Range("D40").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalanceTests.Daily.LeaveBalance_48EmployeeCreation_83_100"
$ws.Range("D41").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalanceTests.Daily.LeaveBalance_48EmployeeCreation_101_123"

# Add new rows 51-55: the new "Hourly" leave-balance test run.
$ws.Range("A51").Value = "56"
$ws.Range("B51").Value = "LeaveBalance_Hourly_1_25"
$ws.Range("C51").Value = "LeaveBalance_Hourly_1_25"
$ws.Range("D51").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalanceTests.Hourly.LeaveBalance_48EmployeeCreation_1_25"
$ws.Range("E51").Value = "Accural//LeaveBalanceTests.xlsx"
$ws.Range("F51").Value = "dummySheet"
$ws.Range("G51").Value = "All"

$ws.Range("A52").Value = "57"
$ws.Range("B52").Value = "LeaveBalance_Hourly_26_41"
$ws.Range("C52").Value = "LeaveBalance_Hourly_26_41"
$ws.Range("D52").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalanceTests.Hourly.LeaveBalance_48EmployeeCreation_26_41"
$ws.Range("E52").Value = "Accural//LeaveBalanceTests.xlsx"
$ws.Range("F52").Value = "dummySheet"
$ws.Range("G52").Value = "All"

$ws.Range("A53").Value = "58"
$ws.Range("B53").Value = "LeaveBalance_Hourly_26_41"
$ws.Range("C53").Value = "LeaveBalance_Hourly_26_41"
$ws.Range("D53").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalanceTests.Hourly.LeaveBalance_48EmployeeCreation_42_82"
$ws.Range("E53").Value = "Accural//LeaveBalanceTests.xlsx"
$ws.Range("F53").Value = "dummySheet"
$ws.Range("G53").Value = "All"

$ws.Range("A54").Value = "59"
$ws.Range("B54").Value = "LeaveBalance_Hourly_83_100"
$ws.Range("C54").Value = "LeaveBalance_Hourly_83_100"
$ws.Range("D54").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalanceTests.Hourly.LeaveBalance_48EmployeeCreation_83_100"
$ws.Range("E54").Value = "Accural//LeaveBalanceTests.xlsx"
$ws.Range("F54").Value = "dummySheet"
$ws.Range("G54").Value = "All"

$ws.Range("A55").Value = "60"
$ws.Range("B55").Value = "LeaveBalance_Hourly_101_123"
$ws.Range("C55").Value = "LeaveBalance_Hourly_101_123"
$ws.Range("D55").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalanceTests.Hourly.LeaveBalance_48EmployeeCreation_101_123"
$ws.Range("E55").Value = "Accural//LeaveBalanceTests.xlsx"
$ws.Range("F55").Value = "dummySheet"
$ws.Range("G55").Value = "All"

# Scroll/selection state, matching the new extent of populated data.
$ws.Range("C51:C55").Select()
$ws.Application.ActiveWindow.ScrollRow = 24
